# Natmi following Dr Hou advice
# Update the Calr-Tshr LR-pairs sheet: recompute existing sending/target
# cluster combinations (rows 2-7) with the new per-gene stats, and add the
# three additional combinations involving the "sCs" sending cluster
# (rows 8-10), extending the sheet from A1:T7 to A1:T10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row (columns A..T) for sheet rows 2..10.
$rowsData = @(
    @("ECs", "Calr", "Tshr", "ECs", 3, 1, 99.883077, 299.649231, 0.3917580942718377, 0.3917580942718377, 3, 1, 0.960636, 2.881908, 0.2124324572954377, 0.2124324572954377, 95.95127955697201, 863.561516012748, 0.0832221346315442, 0.08322213463154422),
    @("ECs", "Calr", "Tshr", "FAPs", 3, 1, 99.883077, 299.649231, 0.3917580942718377, 0.3917580942718377, 2, 0.6666666666666666, 0.5285266666666667, 1.58558, 0.116876963330717, 0.116876963330717, 52.79086974322, 475.11782768898, 0.04578749641872116, 0.04578749641872117),
    @("ECs", "Calr", "Tshr", "sCs", 3, 1, 99.883077, 299.649231, 0.3917580942718377, 0.3917580942718377, 3, 1, 3.032914666666667, 9.098744, 0.6706905793738452, 0.6706905793738454, 302.936849185096, 2726.431642665864, 0.2627484632215723, 0.2627484632215724),
    @("FAPs", "Calr", "Tshr", "ECs", 3, 1, 124.0161413333333, 372.048424, 0.4864119993789693, 0.4864119993789694, 3, 1, 0.960636, 2.881908, 0.2124324572954377, 0.2124324572954377, 119.134369945888, 1072.209329512992, 0.1033296962860613, 0.1033296962860614),
    @("FAPs", "Calr", "Tshr", "FAPs", 3, 1, 124.0161413333333, 372.048424, 0.4864119993789693, 0.4864119993789694, 2, 0.6666666666666666, 0.5285266666666667, 1.58558, 0.116876963330717, 0.116876963330717, 65.54583779176889, 589.91254012592, 0.05685035741503655, 0.05685035741503656),
    @("FAPs", "Calr", "Tshr", "sCs", 3, 1, 124.0161413333333, 372.048424, 0.4864119993789693, 0.4864119993789694, 3, 1, 3.032914666666667, 9.098744, 0.6706905793738452, 0.6706905793738454, 376.1303739532729, 3385.173365579456, 0.3262319456778713, 0.3262319456778715),
    @("sCs", "Calr", "Tshr", "ECs", 3, 1, 31.06188766666667, 93.18566300000001, 0.121829906349193, 0.121829906349193, 3, 1, 0.960636, 2.881908, 0.2124324572954377, 0.2124324572954377, 29.839167520556, 268.552507685004, 0.02588062637783211, 0.02588062637783212),
    @("sCs", "Calr", "Tshr", "FAPs", 3, 1, 31.06188766666667, 93.18566300000001, 0.121829906349193, 0.121829906349193, 2, 0.6666666666666666, 0.5285266666666667, 1.58558, 0.116876963330717, 0.116876963330717, 16.41703594883778, 147.75332353954, 0.01423910949695932, 0.01423910949695932),
    @("sCs", "Calr", "Tshr", "sCs", 3, 1, 31.06188766666667, 93.18566300000001, 0.121829906349193, 0.121829906349193, 3, 1, 3.032914666666667, 9.098744, 0.6706905793738452, 0.6706905793738454, 94.20805467858578, 847.872492107272, 0.08171017047440154, 0.08171017047440159)
)

$rowIndex = 2
foreach ($rowData in $rowsData) {
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
    $rowIndex++
}
